$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "IMF (20%) - Sales" / "IMF (20%) - Sales + Emp" metric pair (stateless
# entities outside the US, taxed at 20%) and drop the old "OECD (20%)" metric pair.
# Layout before: group, M_ETR, M_PL, GFA-Sales, GFA-Sales+Emp, IMF-Sales, IMF-Sales+Emp,
#                OECD(20%)-Sales, OECD(20%)-Sales+Emp, OECD-Sales, OECD-Sales+Emp
# Layout after:  group, M_ETR, M_PL, GFA-Sales, GFA-Sales+Emp, IMF(20%)-Sales, IMF(20%)-Sales+Emp,
#                IMF-Sales, IMF-Sales+Emp, OECD-Sales, OECD-Sales+Emp

# Insert 2 new columns before column F to make room for the new "IMF (20%)" metrics
$ws.Range("F1:G1").EntireColumn.Insert()

# Remove the old "OECD (20%)" columns, which have now been shifted to J:K
$ws.Range("J1:K1").EntireColumn.Delete()

# Set the new column headers
$ws.Range("F1").Value2 = "IMF (20%) - Sales"
$ws.Range("G1").Value2 = "IMF (20%) - Sales + Emp"

# Populate the new "IMF (20%)" data values
$ws.Range("F2").Value2 = -0.0366156072057573
$ws.Range("G2").Value2 = 0.03480246533491808
$ws.Range("F3").Value2 = 0.3454564909433036
$ws.Range("G3").Value2 = 0.311301961979455
$ws.Range("F4").Value2 = 0.8198877780764348
$ws.Range("G4").Value2 = 1.864423691136007
$ws.Range("F5").Value2 = -0.4219475632461457
$ws.Range("G5").Value2 = -0.3017272004193731
$ws.Range("F6").Value2 = 1.508446806846264
$ws.Range("G6").Value2 = 1.052810667356595
$ws.Range("F7").Value2 = -8.108193661997705
$ws.Range("G7").Value2 = 15.86246084502586
$ws.Range("F8").Value2 = 0.50468493164738
$ws.Range("G8").Value2 = 0.6871155749233844
